$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Re-order the sheet tabs: "Demo_TotalOrders" moves in front of
#    "Demo_MultipleProducts" (so the tab order becomes ... Demo_CreateAddress,
#    Demo_TotalOrders, Demo_MultipleProducts, Demo_ApplyDiscount, ...).
#    Excel keeps each sheet's own identity/content attached to its *name*
#    when it is moved, so the data that used to live under the
#    "Demo_MultipleProducts" tab now ends up addressed via that same name.
# ---------------------------------------------------------------------------
$totalOrders = $wb.Worksheets.Item("Demo_TotalOrders")
$multipleProducts = $wb.Worksheets.Item("Demo_MultipleProducts")
$totalOrders.Move($multipleProducts)

# ---------------------------------------------------------------------------
# 2. Populate the "Demo_MultipleProducts" worksheet with its test-data grid.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Demo_MultipleProducts")
$headerSource = $wb.Worksheets.Item("HRM_Login")

$ws.Range("A1").Value = "TestCase ID"
$ws.Range("B1").Value = "TestCase Name"
$ws.Range("C1").Value = "Email"
$ws.Range("D1").Value = "Password"
$ws.Range("E1").Value = "Product1"
$ws.Range("F1").Value = "Produ2"
$ws.Range("G1").Value = "Product3"
$ws.Range("H1").Value = "Product4"
$ws.Range("I1").Value = "Product5"
$ws.Range("J1").Value = "OrderNumber"

# Re-use the existing bold/yellow-fill header style (same one used on every
# other sheet's first row) instead of fabricating a brand-new style entry.
$headerSource.Range("A1").Copy() | Out-Null
$ws.Range("A1:J1").PasteSpecial(-4122) | Out-Null

# Order numbers kept as text. Typing a leading apostrophe mirrors how a user
# would force Excel to store a number-looking value as text (and is what
# produces the "quote prefix" cell style Excel applies automatically).
$ws.Range("J2").Value = "'15623987"
$ws.Range("J3").Value = "'423131"
$ws.Range("J4").Value = "'99823987"
# The last one keeps the plain/default cell style (no quote-prefix marker).
$ws.Range("J4").ClearFormats()

$ws.Columns.Item(1).ColumnWidth = 9.94
$ws.Columns.Item(2).ColumnWidth = 13.5
$ws.Columns.Item(10).ColumnWidth = 12.39

$ws.Range("J3").Select() | Out-Null

# ---------------------------------------------------------------------------
# 3. Selection / active-cell bookkeeping on a handful of other sheets.
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("HRM_Login").Range("B3").Select() | Out-Null
$wb.Worksheets.Item("HRM_AddUser").Range("B3").Select() | Out-Null
$wb.Worksheets.Item("HRM_Performance").Range("N39").Select() | Out-Null

# Leave the selection on Demo_MultipleProducts as the last active state,
# matching its tabSelected="1" sheet view.
$ws.Activate() | Out-Null
$ws.Range("J3").Select() | Out-Null
